$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix column layout: the original "A:B" merged column range (a latent bug)
# is split here so that column A keeps width 30.7109375/style 1 and column B
# keeps its own width 60.7109375/style 2 as two independent <col> entries.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# --- Insert a new row at position 13 (pushes old rows 13-23 down to 14-24,
# carrying their content, styles and row heights along automatically).
$ws.Rows.Item(13).Insert()

# The insert leaves a stray formatted-but-empty cell at A13; the target layout
# has no A-column entry on row 13 at all, so fully clear it.
$ws.Range("A13").Clear()

# New cells B13/C13 did not exist before the insert, so they do not yet carry
# the column's intended "wrap text" style. Copy formatting from existing,
# correctly-styled cells in the same columns before setting their values.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Populate the corrected / new text content.
$ws.Range("B10").Value = 'Fornecer para o aluno habilidades básicas no desenvolvimento de modelos, experimentação, interpretação e análise dos resultados fornecidos através de sessões de simulação, bem como capacitá-lo na escolha da ferramenta adequada para determinadas classes de problemas.'
$ws.Range("C10").Value = 'Fornecer para o aluno habilidades básicas no desenvolvimento de modelos, experimentação, interpretação e análise dos resultados fornecidos através de sessões de simulação, bem como capacitá-lo na escolha da ferramenta adequada para determinadas classes de problemas.'
$ws.Range("B13").Value = '5840917 - Fabricio Maciel Gomes'
$ws.Range("C13").Value = '5840917 - Fabricio Maciel Gomes'
$ws.Range("B14").Value = 'Teoria dos Sistemas, Processo de Modelagem de Sistemas, Modelagem para Simulação, Linguagens de Simulação, Aspectos Estatísticos da Simulação de Sistemas, Validação de Modelos.'
$ws.Range("C14").Value = 'Teoria dos Sistemas, Processo de Modelagem de Sistemas, Modelagem para Simulação, Linguagens de Simulação, Aspectos Estatísticos da Simulação de Sistemas, Validação de Modelos.'
$ws.Range("B16").Value = '1. Teoria dos Sistemas. 1.1. Conceitos Básicos sobre Teoria dos Sistemas;1.2. Abordagem Sistêmica; 1.3. Classificação dos Sistemas; 1.4. Ciclo de Vida dos Sistemas; 2. Processo de Modelagem de Sistemas. 2.1. Modelo Formal; 2.2. Modelo Computacional; 2.3. Teoria dos Modelos; 3. Modelagem para Simulação 3.1. Modelagem Discreta; 3.2. Modelagem Orientada a Eventos; 4. Linguagens de Simulação. 5. Aspectos Estatísticos da Simulação de Sistemas. 5.1. Geradores de Números Aleatórios; 5.2. Geração de Variáveis Aleatórias; 5.3. Inferência Estatística; 5.4. Problemas Estatísticos Relacionados com Simulação; 6. Validação de Modelos.'
$ws.Range("C16").Value = '1. Teoria dos Sistemas. 1.1. Conceitos Básicos sobre Teoria dos Sistemas;1.2. Abordagem Sistêmica; 1.3. Classificação dos Sistemas; 1.4. Ciclo de Vida dos Sistemas; 2. Processo de Modelagem de Sistemas. 2.1. Modelo Formal; 2.2. Modelo Computacional; 2.3. Teoria dos Modelos; 3. Modelagem para Simulação 3.1. Modelagem Discreta; 3.2. Modelagem Orientada a Eventos; 4. Linguagens de Simulação. 5. Aspectos Estatísticos da Simulação de Sistemas. 5.1. Geradores de Números Aleatórios; 5.2. Geração de Variáveis Aleatórias; 5.3. Inferência Estatística; 5.4. Problemas Estatísticos Relacionados com Simulação; 6. Validação de Modelos.'
$ws.Range("B19").Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Range("C19").Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Range("B20").Value = 'NF≥ 5,0.'
$ws.Range("C20").Value = 'NF≥ 5,0.'
$ws.Range("B21").Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'
$ws.Range("C21").Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'
$ws.Range("B22").Value = '1. BANKS, J., CARSON, J.S., “Discrete-Event System Simulation.”, Printice Hall, 2000.2. CHWIF, L., MEDINA, A.C. e col. “Introdução ao Simul8: um guia prático”, Livro Eletrônico, 1ª ed., 2015.3. Fishman, G.S. “Discrete-Event Simulation: Modeling, Programming, and Analysis”, Springer-Verlag, 2001.'
$ws.Range("C22").Value = '1. BANKS, J., CARSON, J.S., “Discrete-Event System Simulation.”, Printice Hall, 2000.2. CHWIF, L., MEDINA, A.C. e col. “Introdução ao Simul8: um guia prático”, Livro Eletrônico, 1ª ed., 2015.3. Fishman, G.S. “Discrete-Event Simulation: Modeling, Programming, and Analysis”, Springer-Verlag, 2001.'

Write-Output "edit complete"
